$wb = $excel.ActiveWorkbook

# Row 46 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I46").Value = 7832
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 23496
$ws.Range("L46").Value = 22500
$ws.Range("M46").Value = -23377
$ws.Range("N46").Value = -22738

# Row 60 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I60").Value = 7832
$ws.Range("J60").Value = 7500
$ws.Range("K60").Value = 23496
$ws.Range("L60").Value = 22500
$ws.Range("M60").Value = -23012
$ws.Range("N60").Value = -23468

# Row 88 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4220.5
$ws.Range("I88").Value = 4218.5
$ws.Range("K88").Value = 4218.5
$ws.Range("M88").Value = -3812.5

# Row 91 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4220.5
$ws.Range("I91").Value = 4218.5
$ws.Range("K91").Value = 4218.5
$ws.Range("M91").Value = -2814.5

# Row 116 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 15924.5
$ws.Range("J116").Value = 7850
$ws.Range("L116").Value = 7850
$ws.Range("N116").Value = -14734

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2644.48
$ws.Range("I138").Value = 1108.069
$ws.Range("J138").Value = 3272.028
$ws.Range("K138").Value = 3324.207
$ws.Range("L138").Value = 9816.083999999999
$ws.Range("M138").Value = 1815.793
$ws.Range("N138").Value = -20096.084

# Row 141 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 15694.667
$ws.Range("I141").Value = 866.24
$ws.Range("K141").Value = 2598.72
$ws.Range("M141").Value = 2581.28

# Row 32 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14559011
$ws.Range("I32").Value = 14707632
$ws.Range("J32").Value = 12825095
$ws.Range("K32").Value = 14707632
$ws.Range("L32").Value = 12825095
$ws.Range("M32").Value = -14707345
$ws.Range("N32").Value = -12825669

# Row 61 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2228.02
$ws.Range("I61").Value = 2103.8372
$ws.Range("K61").Value = 2103.8372
$ws.Range("M61").Value = -1891.8372

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2365
$ws.Range("I74").Value = 2266.423
$ws.Range("K74").Value = 2266.423
$ws.Range("M74").Value = -1392.423

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2365
$ws.Range("I77").Value = 2266.423
$ws.Range("K77").Value = 11332.115
$ws.Range("M77").Value = -6964.114999999998

# Row 110 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1875.9375
$ws.Range("I110").Value = 1644
$ws.Range("J110").Value = 3499.5
$ws.Range("K110").Value = 1644
$ws.Range("L110").Value = 3499.5
$ws.Range("M110").Value = 401
$ws.Range("N110").Value = -7589.5

# Row 136 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2228.02
$ws.Range("I136").Value = 2103.8372
$ws.Range("K136").Value = 6311.5116
$ws.Range("M136").Value = -3761.5116

# Row 105 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2055.739
$ws.Range("I105").Value = 1569.7646
$ws.Range("K105").Value = 1569.7646
$ws.Range("M105").Value = 177.2354

# Row 134 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1788118
$ws.Range("I134").Value = 2042593.8
$ws.Range("J134").Value = 6787.6
$ws.Range("K134").Value = 6127781.4
$ws.Range("L134").Value = 20362.8
$ws.Range("M134").Value = -6125246.4
$ws.Range("N134").Value = -25432.8

# Row 31 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4170.591
$ws.Range("I31").Value = 3286.7273
$ws.Range("J31").Value = 5054.4546
$ws.Range("K31").Value = 3286.7273
$ws.Range("L31").Value = 5054.4546
$ws.Range("M31").Value = -2991.7273
$ws.Range("N31").Value = -5644.4546

# Row 34 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4170.591
$ws.Range("I34").Value = 3286.7273
$ws.Range("J34").Value = 5054.4546
$ws.Range("K34").Value = 3286.7273
$ws.Range("L34").Value = 5054.4546
$ws.Range("M34").Value = -3084.7273
$ws.Range("N34").Value = -5458.4546

# Row 134 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2520.093
$ws.Range("I134").Value = 2300.1765
$ws.Range("J134").Value = 3350.889
$ws.Range("K134").Value = 6900.529500000001
$ws.Range("L134").Value = 10052.667
$ws.Range("M134").Value = -4365.529500000001
$ws.Range("N134").Value = -15122.667

# Row 39 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4949.6
$ws.Range("J39").Value = 4949.6
$ws.Range("L39").Value = 14848.8
$ws.Range("N39").Value = -15436.8

# Row 40 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 138.33333
$ws.Range("I40").Value = 23.636364
$ws.Range("K40").Value = 94.545456
$ws.Range("M40").Value = -25.545456

# Row 50 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 942.9
$ws.Range("J50").Value = 1433
$ws.Range("L50").Value = 4299
$ws.Range("N50").Value = -5261

# Row 53 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 942.9
$ws.Range("J53").Value = 1433
$ws.Range("L53").Value = 4299
$ws.Range("N53").Value = -5261

# Row 55 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1988.909
$ws.Range("J55").Value = 2982.1667
$ws.Range("L55").Value = 8946.500100000001
$ws.Range("N55").Value = -9300.500100000001

# Row 80 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5449.5
$ws.Range("J80").Value = 5900
$ws.Range("L80").Value = 17700
$ws.Range("N80").Value = -19572

# Row 83 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5449.5
$ws.Range("J83").Value = 5900
$ws.Range("L83").Value = 53100
$ws.Range("N83").Value = -62460

# Row 132 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1330.5
$ws.Range("I132").Value = 1364
$ws.Range("J132").Value = 1297
$ws.Range("K132").Value = 12276
$ws.Range("L132").Value = 11673
$ws.Range("M132").Value = -9746
$ws.Range("N132").Value = -16733

# Row 139 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2824.35
$ws.Range("I139").Value = 2571.5454
$ws.Range("J139").Value = 3133.3333
$ws.Range("K139").Value = 7714.6362
$ws.Range("L139").Value = 9399.999899999999
$ws.Range("M139").Value = -2574.6362
$ws.Range("N139").Value = -19679.9999

# Row 96 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 64946.25
$ws.Range("J96").Value = 64946.25
$ws.Range("L96").Value = 64946.25
$ws.Range("N96").Value = -70438.25

# Row 126 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2842.4167
$ws.Range("I126").Value = 2902.25
$ws.Range("J126").Value = 2812.5
$ws.Range("K126").Value = 8706.75
$ws.Range("L126").Value = 8437.5
$ws.Range("M126").Value = -6236.75
$ws.Range("N126").Value = -13377.5

# Row 132 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3071.1191
$ws.Range("I132").Value = 3016
$ws.Range("J132").Value = 3594.75
$ws.Range("K132").Value = 9048
$ws.Range("L132").Value = 10784.25
$ws.Range("M132").Value = -6518
$ws.Range("N132").Value = -15844.25

# Row 42 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 43 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 500005500
$ws.Range("I43").Value = 11012
$ws.Range("K43").Value = 11012
$ws.Range("M43").Value = -10819

# Row 49 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 61 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3773.7896
$ws.Range("I61").Value = 1761.2222
$ws.Range("K61").Value = 1761.2222
$ws.Range("M61").Value = -1559.2222

# Row 113 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3773.7896
$ws.Range("I113").Value = 1761.2222
$ws.Range("K113").Value = 1761.2222
$ws.Range("M113").Value = 408.7778000000001

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 391099.22
$ws.Range("I132").Value = 432524.16
$ws.Range("J132").Value = 4466.6665
$ws.Range("K132").Value = 1297572.48
$ws.Range("L132").Value = 13399.9995
$ws.Range("M132").Value = -1295042.48
$ws.Range("N132").Value = -18459.9995

# Row 136 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2258.282
$ws.Range("I136").Value = 1702.8695
$ws.Range("K136").Value = 5108.6085
$ws.Range("M136").Value = -2558.6085

# Row 117 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 49744.75
$ws.Range("J117").Value = 49744.75
$ws.Range("L117").Value = 49744.75
$ws.Range("N117").Value = -58922.75

# Row 131 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 83281.57000000001
$ws.Range("J131").Value = 83281.57000000001
$ws.Range("L131").Value = 83281.57000000001
$ws.Range("N131").Value = -93361.57000000001

# Row 132 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 40737.926
$ws.Range("I132").Value = 48977.91
$ws.Range("K132").Value = 146933.73
$ws.Range("M132").Value = -144403.73

# Row 136 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 30024.5
$ws.Range("I136").Value = 1922.6072
$ws.Range("K136").Value = 5767.821599999999
$ws.Range("M136").Value = -3217.821599999999
